# GlobalConfig.xlsx balance-tuning update (晶石/装备 经验道具 & 升级参数)
#
# - equipExpItems (装备经验道具id) value:
#     150,151,152,153,154 -> 400,401,402,403,404
# - crystalSwallowExpLoss (晶石吞噬经验折损率):
#     type int32 -> number, value 8000 -> 0.8
# - crystalExpItems (晶石经验道具id) value:
#     200,201,202,203,204 -> 600,601,602,603,604
# - crystalLevelupQualityRatio (晶石升级品质系数) value:
#     0.7,0.85,1,1.15,1.25 -> 0.7,0.85,1,1.15,1.3
# - crystalLevelupRandRatio (晶石升级副属性随机区间系数) value:
#     0.9,1.1 -> 0.7,1.1
# - crystalLevelupAssistantNumber (晶石副属性随机到相同属性的次数上限) value:
#     3 -> 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G21").Value = "400,401,402,403,404"

$ws.Range("F32").Value = "number"
$ws.Range("G32").Value = "0.8"

$ws.Range("G34").Value = "600,601,602,603,604"

$ws.Range("G36").Value = "0.7,0.85,1,1.15,1.3"

$ws.Range("G37").Value = "0.7,1.1"

$ws.Range("G40").Value = "6"

# Move the cursor/selection the same way the original author left it.
$ws.Range("F42").Select()
